# Orders.xlsx update: "filemanager read csv settings"
#
# 1) On sheet "Sheet": rows 61-80 get a single-space placeholder in column D
#    (an empty employer-name slot), and rows 125-133 get the employer name
#    "Ludmyla Stupnytska" filled in to column D (replacing the placeholder
#    space).
# 2) On sheet "Employers": the "Maryna Mostishko" employer column (G:H) is
#    removed entirely, and nine new order rows (22-30) are appended under
#    the "Ludmyla Stupnytska" column (C:D).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "Sheet"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sheet")

for ($r = 61; $r -le 80; $r++) {
    $ws.Cells.Item($r, 4).Value = " "
}

for ($r = 125; $r -le 133; $r++) {
    $ws.Cells.Item($r, 4).Value = "Ludmyla Stupnytska"
}

# ---------------------------------------------------------------------
# Sheet: "Employers"
# ---------------------------------------------------------------------
$es = $wb.Worksheets.Item("Employers")

# Drop the "Maryna Mostishko" employer (columns G:H, header + 20 data rows).
$es.Range("G1:H21").EntireColumn.Delete()

# Append the new "Ludmyla Stupnytska" order codes in columns C:D, rows 22-30.
# Force text formatting so the dates stay literal strings ("2020-09-30")
# instead of being auto-parsed into date serials, matching the rest of
# the sheet.
$es.Range("D22:D30").NumberFormat = "@"

$newCodes = @(
    "GRANNL1531234",
    "GRANNL1531229",
    "GRANNL1531223",
    "GRANNL1531222",
    "GRANNL1531220",
    "GRANNL1531217",
    "GRANNL1531216",
    "GRANNL1531214",
    "GRANNL1531212"
)

$r = 22
foreach ($code in $newCodes) {
    $es.Cells.Item($r, 3).Value = $code
    $es.Cells.Item($r, 4).Value = "2020-09-30"
    $r = $r + 1
}
